$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 78, shifting the existing rows 78:147 down to 79:148.
$ws.Rows(78).Insert()

# Populate the newly inserted row 78 with the new record.
$ws.Range("A78").Value = 9
$ws.Range("B78").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C78").Value = "Metropolitana"
$ws.Range("D78").Value = 44897
$ws.Range("E78").Value = 13
$ws.Range("F78").Value = 100112022
$ws.Range("G78").Value = "Arveja Verde"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 90
$ws.Range("K78").Value = 20000
$ws.Range("L78").Value = 21000
$ws.Range("M78").Value = 20500
$ws.Range("N78").Value = "$/saco 25 kilos"
$ws.Range("O78").Value = "Región del Maule"
$ws.Range("P78").Value = 820
$ws.Range("Q78").Value = 25
$ws.Range("R78").Value = "Hortaliza"
